$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.575.50'
$ws.Range("E2").Value = '  -2.94%  '
$ws.Range("D3").Value = '2.712.87'
$ws.Range("E3").Value = '  -6.58%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = '''502.61'
$ws.Range("E5").Value = '  -4.71%  '
$ws.Range("D6").Value = '''139.45'
$ws.Range("E6").Value = '  -1.49%  '
$ws.Range("D7").Value = '''0.996'
$ws.Range("E7").Value = '  -0.46%  '
$ws.Range("D8").Value = '''0.528'
$ws.Range("E8").Value = '  -3.97%  '
$ws.Range("D9").Value = '2.729.37'
$ws.Range("E9").Value = '  -6.08%  '
$ws.Range("D10").Value = '''6.02'
$ws.Range("E10").Value = '  +2.85%  '
$ws.Range("D11").Value = '''0.103'
$ws.Range("E11").Value = '  -4.17%  '
$ws.Range("D12").Value = '''0.343'
$ws.Range("E12").Value = '  -2.48%  '
$ws.Range("E13").Value = '  +1.24%  '
$ws.Range("D14").Value = '3.170.83'
$ws.Range("E14").Value = '  -7.05%  '
$ws.Range("D15").Value = '58.508.13'
$ws.Range("E15").Value = '  -3.28%  '
$ws.Range("D16").Value = '''21.39'
$ws.Range("E16").Value = '  -5.32%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '''0.0000134'
$ws.Range("E17").Value = '  -4.43%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '2.693.52'
$ws.Range("E18").Value = '  -7.30%  '
$ws.Range("D19").Value = '''4.66'
$ws.Range("E19").Value = '  -5.03%  '
$ws.Range("D20").Value = '''10.83'
$ws.Range("E20").Value = '  -5.83%  '
$ws.Range("D21").Value = '''338.71'
$ws.Range("E21").Value = '  -5.96%  '
$ws.Range("D22").Value = '''6.16'
$ws.Range("E22").Value = '  -6.50%  '
$ws.Range("D23").Value = '''1.00'
$ws.Range("E23").Value = '  +0.27%  '
$ws.Range("D24").Value = '''5.60'
$ws.Range("E24").Value = '  -0.71%  '
$ws.Range("D25").Value = '''62.43'
$ws.Range("E25").Value = '  -1.59%  '
$ws.Range("D26").Value = '''0.420'
$ws.Range("E26").Value = '  -5.90%  '
$ws.Range("E27").Value = '  -3.57%  '
$ws.Range("D28").Value = '''0.994'
$ws.Range("E28").Value = '  -0.52%  '
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").Value = '0.0₃0821'
$ws.Range("E29").Value = '  -3.28%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").Value = '''7.39'
$ws.Range("E30").Value = '  -3.75%  '
$ws.Range("D31").Value = '''0.997'
$ws.Range("E31").Value = '  -0.25%  '
$ws.Range("D32").Value = '''1.59'
$ws.Range("E32").Value = '  -4.13%  '
$ws.Range("D33").Value = '''18.91'
$ws.Range("E33").Value = '  -3.50%  '
$ws.Range("D34").Value = '''147.70'
$ws.Range("E34").Value = '  -3.18%  '
$ws.Range("D35").Value = '''4.13'
$ws.Range("E35").Value = '  -3.52%  '
$ws.Range("D36").Value = '''5.30'
$ws.Range("E36").Value = '  -3.92%  '
$ws.Range("D37").Value = '''0.927'
$ws.Range("E37").Value = '  -6.26%  '
$ws.Range("D38").Value = '''1.12'
$ws.Range("E38").Value = '  -5.28%  '
$ws.Range("D39").Value = '''36.07'
$ws.Range("E39").Value = '  -4.69%  '
$ws.Range("E40").Value = '  -5.47%  '
$ws.Range("D41").Value = '2.160.23'
$ws.Range("E41").Value = '  -7.45%  '
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").Value = '''3.48'
$ws.Range("E42").Value = '  -4.61%  '
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").Value = '''0.995'
$ws.Range("E43").Value = '  -0.39%  '
$ws.Range("D44").Value = '''0.0550'
$ws.Range("E44").Value = '  -3.01%  '
$ws.Range("D45").Value = '''0.597'
$ws.Range("E45").Value = '  -6.99%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '''18.81'
$ws.Range("E46").Value = '  -9.19%  '
$ws.Range("B47").Value = 'WhiteBITCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D47").Value = '''10.34'
$ws.Range("E47").Value = '  -0.13%  '
$ws.Range("D48").Value = '''0.0224'
$ws.Range("E48").Value = '  -3.42%  '
$ws.Range("D49").Value = '''4.58'
$ws.Range("E49").Value = '  -4.79%  '
$ws.Range("D50").Value = '''0.0880'
$ws.Range("E50").Value = '  -4.58%  '
$ws.Range("D51").Value = '''17.70'
$ws.Range("E51").Value = '  -2.63%  '
